# Updates the cryptos price/ranking table to the values scraped on
# Fri Dec 30 19:21:14 UTC 2022.
#
# All "Price" column (column D) values are stored as text in the
# workbook, so we force a text NumberFormat before assigning them and
# then reset the cell style back to Normal so no visible formatting
# change is introduced; this preserves exact text such as trailing
# zeros (e.g. "6.540", "0.1900") instead of Excel coercing the value to
# a numeric type.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "245.25"
$ws.Cells.Item(2, 4).Style = "Normal"

$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = "25.14"
$ws.Cells.Item(3, 4).Style = "Normal"

$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = "5.089"
$ws.Cells.Item(4, 4).Style = "Normal"

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "6.540"
$ws.Cells.Item(6, 4).Style = "Normal"

$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "3.004"
$ws.Cells.Item(7, 4).Style = "Normal"

$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "0.8094"
$ws.Cells.Item(8, 4).Style = "Normal"

$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.8383"
$ws.Cells.Item(9, 4).Style = "Normal"

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "0.1337"
$ws.Cells.Item(10, 4).Style = "Normal"

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.06940"
$ws.Cells.Item(11, 4).Style = "Normal"

$ws.Cells.Item(12, 2).Value = "LiechtensteinCryptoassetsExchange"

$ws.Cells.Item(12, 3).Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"

$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "0.03232"
$ws.Cells.Item(12, 4).Style = "Normal"

$ws.Cells.Item(12, 5).Value = "11LiechtensteinCryptoassetsExchangeLCX"

$ws.Cells.Item(13, 2).Value = "BitrueCoin"

$ws.Cells.Item(13, 3).Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"

$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "0.02821"
$ws.Cells.Item(13, 4).Style = "Normal"

$ws.Cells.Item(13, 5).Value = "12BitrueCoinBTR"

$ws.Cells.Item(14, 2).Value = "BitMartToken"

$ws.Cells.Item(14, 3).Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"

$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "0.09412"
$ws.Cells.Item(14, 4).Style = "Normal"

$ws.Cells.Item(14, 5).Value = "13BitMartTokenBMX"

$ws.Cells.Item(15, 2).Value = "BitForexToken"

$ws.Cells.Item(15, 3).Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"

$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "0.001509"
$ws.Cells.Item(15, 4).Style = "Normal"

$ws.Cells.Item(15, 5).Value = "14BitForexTokenBF"

$ws.Cells.Item(16, 2).Value = "One"

$ws.Cells.Item(16, 3).Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"

$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "0.0005989"
$ws.Cells.Item(16, 4).Style = "Normal"

$ws.Cells.Item(16, 5).Value = "15OneONE"

$ws.Cells.Item(17, 2).Value = "TigerCash"

$ws.Cells.Item(17, 3).Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"

$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "0.006089"
$ws.Cells.Item(17, 4).Style = "Normal"

$ws.Cells.Item(17, 5).Value = "16TigerCashTCH"

$ws.Cells.Item(18, 2).Value = "LEO"

$ws.Cells.Item(18, 3).Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"

$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "3.500"
$ws.Cells.Item(18, 4).Style = "Normal"

$ws.Cells.Item(18, 5).Value = "17LEOLEO"

$ws.Cells.Item(19, 2).Value = "BTSEToken"

$ws.Cells.Item(19, 3).Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"

$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "2.092"
$ws.Cells.Item(19, 4).Style = "Normal"

$ws.Cells.Item(19, 5).Value = "18BTSETokenBTSE"

$ws.Cells.Item(20, 2).Value = "BitpandaEcosystemToken"

$ws.Cells.Item(20, 3).Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"

$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "0.3182"
$ws.Cells.Item(20, 4).Style = "Normal"

$ws.Cells.Item(20, 5).Value = "19BitpandaEcosystemTokenBEST"

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "3.765"
$ws.Cells.Item(22, 4).Style = "Normal"

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "0.04689"
$ws.Cells.Item(23, 4).Style = "Normal"

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "0.001241"
$ws.Cells.Item(25, 4).Style = "Normal"

$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "0.004527"
$ws.Cells.Item(26, 4).Style = "Normal"

$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "0.00009698"
$ws.Cells.Item(27, 4).Style = "Normal"

$ws.Cells.Item(27, 5).Value = "26NitroExNTX"

$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "0.1347"
$ws.Cells.Item(41, 4).Style = "Normal"

$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "0.006239"
$ws.Cells.Item(42, 4).Style = "Normal"

$ws.Cells.Item(42, 5).Value = "41KickTokenKICKBestin24h"

$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "0.002722"
$ws.Cells.Item(43, 4).Style = "Normal"

$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "0.008054"
$ws.Cells.Item(44, 4).Style = "Normal"

$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "0.00005278"
$ws.Cells.Item(45, 4).Style = "Normal"

$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "0.00000000750"
$ws.Cells.Item(46, 4).Style = "Normal"

$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "0.1900"
$ws.Cells.Item(47, 4).Style = "Normal"

$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "0.002046"
$ws.Cells.Item(48, 4).Style = "Normal"
